# Apply updated balance-sheet figures to the "GNRC" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GNRC")

# Row 4 (Inventory)
$ws.Range("B4").Value = 645000000.0
$ws.Range("C4").Value = 603000000.0
$ws.Range("D4").Value = 533000000.0
$ws.Range("E4").Value = 544000000.0
$ws.Range("F4").Value = 560000000.0

# Row 13 (Accounts Payable)
$ws.Range("B13").Value = 388000000.0
$ws.Range("C13").Value = 330000000.0
$ws.Range("D13").Value = 273000000.0
$ws.Range("E13").Value = 230000000.0
$ws.Range("F13").Value = 267000000.0

# Row 20 (Long Term Tax Liability (Deferred))
$ws.Range("B20").Value = 121000000.0
$ws.Range("C20").Value = 114000000.0
$ws.Range("D20").Value = 102000000.0
$ws.Range("E20").Value = 94000000.0
$ws.Range("F20").Value = 89000000.0

# Row 35 (Net Debt)
$ws.Range("G35").Value = 575981000.0

# Row 36 (Total Debt)
$ws.Range("G36").Value = 898864000.0
